$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header cell C1: "Function" -> "Drawing Type"
$ws.Range("C1").Value = "Drawing Type"

# Update selection to C12 (as shown in the diff's sheetView selection)
$ws.Range("C12").Select()
